$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------------
# New source file being handed off: e1d1d644-a81d-4eb3-9627-5f9d64c9ca9d.md
# Add one new row to each of the three tables (Overview, zh-cn, de-de).
# -------------------------------------------------------------------------

$mdFile        = "e1d1d644-a81d-4eb3-9627-5f9d64c9ca9d.md"
$mdPath        = "e2e\e1d1d644-a81d-4eb3-9627-5f9d64c9ca9d.md"
$githubUrl     = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ccc8e94409f51ffb777cabd56c773bd6ad332c5f/e2e/e1d1d644-a81d-4eb3-9627-5f9d64c9ca9d.md"
$dateFmt       = "yyyy-mm-dd HH:mm:ss"
$hyperlinkRGB  = 15570276   # OLE BGR for #6495ED, matches the workbook's existing HyperLink font color

# ---------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = $mdFile
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $githubUrl, "", "", $mdPath) | Out-Null
$wsOverview.Range("B3").Font.Underline = $true
$wsOverview.Range("B3").Font.Color = $hyperlinkRGB
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").NumberFormat = $dateFmt
$wsOverview.Range("G3").Value = "2016-09-05 16:47:34"

# ---------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $githubUrl, "", "", $mdFile) | Out-Null
$wsZhCn.Range("A3").Font.Underline = $true
$wsZhCn.Range("A3").Font.Color = $hyperlinkRGB
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("G3").Value = "e1d1d644-a81d-4eb3-9627-5f9d64c9ca9d.bcf5fe067bc88e2c8daa1e22cca9ed4c8163c20b.zh-cn.xlf"
$wsZhCn.Range("H3").NumberFormat = $dateFmt
$wsZhCn.Range("H3").Value = "2016-09-05 16:47:30"
$wsZhCn.Range("K3").NumberFormat = $dateFmt
$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("M3").Value = "'True"
$wsZhCn.Range("O3").Value = "'False"

# ---------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $githubUrl, "", "", $mdFile) | Out-Null
$wsDeDe.Range("A3").Font.Underline = $true
$wsDeDe.Range("A3").Font.Color = $hyperlinkRGB
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("G3").Value = "e1d1d644-a81d-4eb3-9627-5f9d64c9ca9d.bcf5fe067bc88e2c8daa1e22cca9ed4c8163c20b.de-de.xlf"
$wsDeDe.Range("H3").NumberFormat = $dateFmt
$wsDeDe.Range("H3").Value = "2016-09-05 16:47:34"
$wsDeDe.Range("K3").NumberFormat = $dateFmt
$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("M3").Value = "'True"
$wsDeDe.Range("O3").Value = "'False"
